$d = $word.ActiveDocument

# Guard text at the very end of the document. Several Range-position-based ops
# below (in particular Bookmarks.Add) are unreliable right at the tail of the
# story; parking harmless text there first keeps every position we compute
# comfortably away from the document end, and we strip it again at the end.
$guard = $d.Content
$guard.Collapse(0)
$guard.InsertAfter("GUARDTEXTGUARDTEXTGUARDTEXT")

# 1) "Принцип 5?." currently sits in three runs split by gramStart/gramEnd
#    proofErr markers; collapse it back into a single plain run.
$r0 = $d.Content
$r0.Find.Execute("Принцип 5?. Психофизиологическое состояние имеет значение", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "Принцип 5?. Психофизиологическое состояние имеет значение", 2) | Out-Null

# 2) Locate the end of "... Иметь связь с картинкой, запахом." -- that's
#    where the new "Действие 8" paragraph needs to start.
$targetRange = $d.Content
$targetRange.Find.Execute("Иметь связь с картинкой, запахом.") | Out-Null
$endBeforeSplit = $targetRange.End

# 3) The _GoBack bookmark currently wraps the end of that paragraph; drop it
#    here and re-create it at the end of the new paragraph below.
$d.Bookmarks.Item("_GoBack").Delete()

# 4) Split so "Действие 8" begins a brand-new paragraph right after it.
$targetRange.InsertParagraphAfter()

# 5) Find the paragraph that was just created (it starts right after the
#    inserted paragraph mark) and give it its text.
$paras = $d.Paragraphs
$newParaIndex = $null
for ($i = 1; $i -le $paras.Count; $i++) {
    if ($paras.Item($i).Range.Start -eq ($endBeforeSplit + 1)) {
        $newParaIndex = $i
    }
}
$newPara = $paras.Item($newParaIndex)
$newPara.Range.InsertAfter("Действие 8. Еще не завршено")

# 6) Re-create the _GoBack bookmark right after the new run, inside the new
#    paragraph. Bookmarks.Add is unreliable when the target position sits
#    immediately before a paragraph mark, so nudge the mark out of the way
#    with a one-character local guard, add the bookmark, then delete the
#    guard character again.
$paras2 = $d.Paragraphs
$newPara2 = $paras2.Item($newParaIndex)
$bmPos = $newPara2.Range.End - 1
$d.Range($bmPos, $bmPos).InsertAfter("X")
$d.Bookmarks.Add("_GoBack", $d.Range($bmPos, $bmPos))
$d.Range($bmPos, $bmPos + 1).Delete()

# 7) Remove the guard text parked at the end of the document in step 0.
$gr = $d.Content
$gr.Find.Execute("GUARDTEXTGUARDTEXTGUARDTEXT") | Out-Null
$gr.Delete()
